# Add a second vehicle (an "electricity" vehicle, e.g. an EV) as row 3 of the
# OperationScenario_Vehicle table, mirroring the unit/label columns of row 2
# but with the numeric parameters zeroed out (capacity, consumption rate,
# charge/discharge efficiency and power all set to 0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "electricity"
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = "Wh"
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = "Wh/km"
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = "W"
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = "W"
$ws.Range("M3").Value = 0

# Grow the table (ListObject) so the new row is part of Table1.
$table = $ws.ListObjects.Item("Table1")
$table.Resize($ws.Range("A1:M4"))

# Matches the saved selection position recorded in the workbook.
$ws.Range("D8").Select()
